$wb = $excel.ActiveWorkbook

function Set-CellValue($ws, $ref, $value) {
    $ws.Range($ref).Value = $value
}

# Sheet: Citywide Totals
$ws = $wb.Worksheets.Item("Citywide Totals")
Set-CellValue $ws "L2" 979
Set-CellValue $ws "L3" 1003
Set-CellValue $ws "D4" 1989
Set-CellValue $ws "E4" 2046
Set-CellValue $ws "F4" 1922
Set-CellValue $ws "L4" 280
Set-CellValue $ws "L5" 65
Set-CellValue $ws "L6" 1033
Set-CellValue $ws "D7" 28180
Set-CellValue $ws "E7" 26052
Set-CellValue $ws "F7" 24115
Set-CellValue $ws "L7" 3360

# Sheet: Austin
$ws = $wb.Worksheets.Item("Austin")
Set-CellValue $ws "L3" 62
Set-CellValue $ws "K4" 103
Set-CellValue $ws "L4" 15
Set-CellValue $ws "L5" 5
Set-CellValue $ws "K7" 1802
Set-CellValue $ws "L7" 202

# Sheet: South Chicago
$ws = $wb.Worksheets.Item("South Chicago")
Set-CellValue $ws "L3" 33
Set-CellValue $ws "L7" 72

# Sheet: Garfield Park
$ws = $wb.Worksheets.Item("Garfield Park")
Set-CellValue $ws "L3" 55
Set-CellValue $ws "K6" 364
Set-CellValue $ws "L6" 54
Set-CellValue $ws "K7" 1150
Set-CellValue $ws "L7" 147

# Sheet: Grand Crossing
$ws = $wb.Worksheets.Item("Grand Crossing")
Set-CellValue $ws "L3" 28
Set-CellValue $ws "L4" 9
Set-CellValue $ws "L7" 117

# Sheet: Woodlawn
$ws = $wb.Worksheets.Item("Woodlawn")
Set-CellValue $ws "L2" 13
Set-CellValue $ws "L3" 26
Set-CellValue $ws "L7" 53

# Sheet: Fuller Park
$ws = $wb.Worksheets.Item("Fuller Park")
Set-CellValue $ws "L6" 7
Set-CellValue $ws "L7" 18

# Sheet: By Neighborhood
$ws = $wb.Worksheets.Item("By Neighborhood")
Set-CellValue $ws "L2" 23
Set-CellValue $ws "L7" 109
Set-CellValue $ws "K8" 1802
Set-CellValue $ws "L8" 202
Set-CellValue $ws "L10" 22
Set-CellValue $ws "L11" 48
Set-CellValue $ws "L21" 10
Set-CellValue $ws "L25" 17
Set-CellValue $ws "L29" 157
Set-CellValue $ws "L30" 18
Set-CellValue $ws "L31" 36
Set-CellValue $ws "K33" 1150
Set-CellValue $ws "L33" 147
Set-CellValue $ws "L36" 60
Set-CellValue $ws "L37" 117
Set-CellValue $ws "K42" 1026
Set-CellValue $ws "L42" 107
Set-CellValue $ws "L44" 23
Set-CellValue $ws "L47" 28
Set-CellValue $ws "L51" 45
Set-CellValue $ws "L54" 73
Set-CellValue $ws "D63" 369
Set-CellValue $ws "E63" 381
Set-CellValue $ws "F63" 207
Set-CellValue $ws "K63" 78
Set-CellValue $ws "L63" 19
Set-CellValue $ws "L67" 123
Set-CellValue $ws "L71" 10
Set-CellValue $ws "L72" 13
Set-CellValue $ws "L75" 14
Set-CellValue $ws "L79" 96
Set-CellValue $ws "L83" 72
Set-CellValue $ws "L84" 30
Set-CellValue $ws "L85" 172
Set-CellValue $ws "L88" 50
Set-CellValue $ws "L89" 46
Set-CellValue $ws "L93" 17
Set-CellValue $ws "L98" 25
Set-CellValue $ws "L99" 53
Set-CellValue $ws "L100" 5
Set-CellValue $ws "D101" 28180
Set-CellValue $ws "E101" 26052
Set-CellValue $ws "F101" 24115
Set-CellValue $ws "L101" 3360

# Sheet: Gage Park
$ws = $wb.Worksheets.Item("Gage Park")
Set-CellValue $ws "L3" 8
Set-CellValue $ws "L7" 36

# Sheet: North Lawndale
$ws = $wb.Worksheets.Item("North Lawndale")
Set-CellValue $ws "L2" 38
Set-CellValue $ws "L5" 5
Set-CellValue $ws "L7" 123

# Sheet: South Deering
$ws = $wb.Worksheets.Item("South Deering")
Set-CellValue $ws "L2" 12
Set-CellValue $ws "L7" 30

# Sheet: Loop
$ws = $wb.Worksheets.Item("Loop")
Set-CellValue $ws "L6" 36
Set-CellValue $ws "L7" 73

# Sheet: Englewood
$ws = $wb.Worksheets.Item("Englewood")
Set-CellValue $ws "L6" 44
Set-CellValue $ws "L7" 157

# Sheet: Irving Park
$ws = $wb.Worksheets.Item("Irving Park")
Set-CellValue $ws "L3" 6
Set-CellValue $ws "L6" 11
Set-CellValue $ws "L7" 23

# Sheet: Humboldt Park
$ws = $wb.Worksheets.Item("Humboldt Park")
Set-CellValue $ws "L2" 28
Set-CellValue $ws "L3" 23
Set-CellValue $ws "K6" 395
Set-CellValue $ws "L6" 46
Set-CellValue $ws "K7" 1026
Set-CellValue $ws "L7" 107

# Sheet: Avondale
$ws = $wb.Worksheets.Item("Avondale")
Set-CellValue $ws "L3" 4
Set-CellValue $ws "L7" 22

# Sheet: Chinatown
$ws = $wb.Worksheets.Item("Chinatown")
Set-CellValue $ws "L3" 4
Set-CellValue $ws "L7" 10

# Sheet: Roseland
$ws = $wb.Worksheets.Item("Roseland")
Set-CellValue $ws "L4" 8
Set-CellValue $ws "L7" 96

# Sheet: Grand Boulevard
$ws = $wb.Worksheets.Item("Grand Boulevard")
Set-CellValue $ws "L3" 14
Set-CellValue $ws "L7" 60

# Sheet: West Lawn
$ws = $wb.Worksheets.Item("West Lawn")
Set-CellValue $ws "L6" 6
Set-CellValue $ws "L7" 17

# Sheet: Wrigleyville
$ws = $wb.Worksheets.Item("Wrigleyville")
Set-CellValue $ws "L3" 2
Set-CellValue $ws "L7" 5

# Sheet: Auburn Gresham
$ws = $wb.Worksheets.Item("Auburn Gresham")
Set-CellValue $ws "L3" 40
Set-CellValue $ws "L7" 109

# Sheet: East Side
$ws = $wb.Worksheets.Item("East Side")
Set-CellValue $ws "L2" 7
Set-CellValue $ws "L7" 17

# Sheet: Kenwood
$ws = $wb.Worksheets.Item("Kenwood")
Set-CellValue $ws "L3" 8
Set-CellValue $ws "L7" 28

# Sheet: Wicker Park
$ws = $wb.Worksheets.Item("Wicker Park")
Set-CellValue $ws "L2" 6
Set-CellValue $ws "L7" 25

# Sheet: Belmont Cragin
$ws = $wb.Worksheets.Item("Belmont Cragin")
Set-CellValue $ws "L6" 15
Set-CellValue $ws "L7" 48

# Sheet: Albany Park
$ws = $wb.Worksheets.Item("Albany Park")
Set-CellValue $ws "L2" 6
Set-CellValue $ws "L7" 23

# Sheet: United Center
$ws = $wb.Worksheets.Item("United Center")
Set-CellValue $ws "L6" 23
Set-CellValue $ws "L7" 50

# Sheet: Uptown
$ws = $wb.Worksheets.Item("Uptown")
Set-CellValue $ws "L4" 8
Set-CellValue $ws "L7" 46

# Sheet: Pullman
$ws = $wb.Worksheets.Item("Pullman")
Set-CellValue $ws "L6" 1
Set-CellValue $ws "L7" 14

# Sheet: Little Italy, UIC
$ws = $wb.Worksheets.Item("Little Italy, UIC")
Set-CellValue $ws "L4" 3
Set-CellValue $ws "L6" 14
Set-CellValue $ws "L7" 45

# Sheet: South Shore
$ws = $wb.Worksheets.Item("South Shore")
Set-CellValue $ws "L3" 74
Set-CellValue $ws "L4" 13
Set-CellValue $ws "L6" 37
Set-CellValue $ws "L7" 172

# Sheet: Oakland
$ws = $wb.Worksheets.Item("Oakland")
Set-CellValue $ws "L6" 4
Set-CellValue $ws "L7" 10

# Sheet: Old Town
$ws = $wb.Worksheets.Item("Old Town")
Set-CellValue $ws "L6" 4
Set-CellValue $ws "L7" 13
